$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1295.746
$ws.Range("I15").Value = 1295.746
$ws.Range("K15").Value = 3887.238
$ws.Range("M15").Value = -3718.238
$ws.Range("H137").Value = 6067.44
$ws.Range("I137").Value = 4917.8125
$ws.Range("K137").Value = 14753.4375
$ws.Range("M137").Value = -12203.4375
$ws.Range("H138").Value = 3306.7327
$ws.Range("I138").Value = 996.381
$ws.Range("J138").Value = 5512.0684
$ws.Range("K138").Value = 2989.143
$ws.Range("L138").Value = 16536.2052
$ws.Range("M138").Value = 2150.857
$ws.Range("N138").Value = -26816.2052
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 52636588
$ws.Range("I2").Value = 3592.5
$ws.Range("J2").Value = 111117700
$ws.Range("K2").Value = 3592.5
$ws.Range("L2").Value = 111117700
$ws.Range("M2").Value = -3479.5
$ws.Range("N2").Value = -111117926
$ws.Range("H32").Value = 1566341.4
$ws.Range("I32").Value = 1605862.5
$ws.Range("K32").Value = 1605862.5
$ws.Range("M32").Value = -1605575.5
$ws.Range("H63").Value = 6001
$ws.Range("I63").Value = 6001
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 6001
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -5315
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 6001
$ws.Range("I66").Value = 6001
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 30005
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -26573
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 36689.375
$ws.Range("I74").Value = 48265.4
$ws.Range("K74").Value = 48265.4
$ws.Range("M74").Value = -47391.4
$ws.Range("H77").Value = 36689.375
$ws.Range("I77").Value = 48265.4
$ws.Range("K77").Value = 241327
$ws.Range("M77").Value = -236959
$ws.Range("H104").Value = 41111
$ws.Range("J104").Value = 41111
$ws.Range("L104").Value = 41111
$ws.Range("N104").Value = -48099
$ws.Range("H116").Value = 52636588
$ws.Range("I116").Value = 3592.5
$ws.Range("J116").Value = 111117700
$ws.Range("K116").Value = 3592.5
$ws.Range("L116").Value = 111117700
$ws.Range("M116").Value = -1298.5
$ws.Range("N116").Value = -111122288
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 52636588
$ws.Range("I3").Value = 3592.5
$ws.Range("J3").Value = 111117700
$ws.Range("K3").Value = 3592.5
$ws.Range("L3").Value = 111117700
$ws.Range("M3").Value = -3478.5
$ws.Range("N3").Value = -111117928
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H82").Value = 16666.334
$ws.Range("I82").Value = 16666.334
$ws.Range("K82").Value = 16666.334
$ws.Range("M82").Value = -16283.334
$ws.Range("H85").Value = 16666.334
$ws.Range("I85").Value = 16666.334
$ws.Range("K85").Value = 16666.334
$ws.Range("M85").Value = -15340.334
$ws.Range("H86").Value = 74252.21000000001
$ws.Range("I86").Value = 112860.336
$ws.Range("J86").Value = 4757.6
$ws.Range("K86").Value = 112860.336
$ws.Range("L86").Value = 4757.6
$ws.Range("M86").Value = -111737.336
$ws.Range("N86").Value = -7003.6
$ws.Range("H88").Value = 64341.5
$ws.Range("J88").Value = 64341.5
$ws.Range("L88").Value = 64341.5
$ws.Range("N88").Value = -65153.5
$ws.Range("H89").Value = 74252.21000000001
$ws.Range("I89").Value = 112860.336
$ws.Range("J89").Value = 4757.6
$ws.Range("K89").Value = 564301.6799999999
$ws.Range("L89").Value = 23788
$ws.Range("M89").Value = -558685.6799999999
$ws.Range("N89").Value = -35020
$ws.Range("H91").Value = 64341.5
$ws.Range("J91").Value = 64341.5
$ws.Range("L91").Value = 64341.5
$ws.Range("N91").Value = -67149.5
$ws.Range("H134").Value = 6255807.5
$ws.Range("I134").Value = 12501420
$ws.Range("J134").Value = 10194.7
$ws.Range("K134").Value = 37504260
$ws.Range("L134").Value = 30584.1
$ws.Range("M134").Value = -37501725
$ws.Range("N134").Value = -35654.10000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7253646.5
$ws.Range("I31").Value = 2280.3333
$ws.Range("J31").Value = 10762372
$ws.Range("K31").Value = 2280.3333
$ws.Range("L31").Value = 10762372
$ws.Range("M31").Value = -1985.3333
$ws.Range("N31").Value = -10762962
$ws.Range("H34").Value = 7253646.5
$ws.Range("I34").Value = 2280.3333
$ws.Range("J34").Value = 10762372
$ws.Range("K34").Value = 2280.3333
$ws.Range("L34").Value = 10762372
$ws.Range("M34").Value = -2078.3333
$ws.Range("N34").Value = -10762776
$ws.Range("H105").Value = 7148763.5
$ws.Range("I105").Value = 14288750
$ws.Range("K105").Value = 14288750
$ws.Range("M105").Value = -14287003
$ws.Range("H132").Value = 14292153
$ws.Range("I132").Value = 2950.1428
$ws.Range("J132").Value = 28581356
$ws.Range("K132").Value = 8850.428400000001
$ws.Range("L132").Value = 85744068
$ws.Range("M132").Value = -6320.428400000001
$ws.Range("N132").Value = -85749128
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 2100
$ws.Range("J28").Value = 2100
$ws.Range("L28").Value = 2100
$ws.Range("N28").Value = -2484
$ws.Range("H29").Value = 3831.3333
$ws.Range("I29").Value = 2999.5
$ws.Range("J29").Value = 5495
$ws.Range("K29").Value = 2999.5
$ws.Range("L29").Value = 5495
$ws.Range("M29").Value = -2709.5
$ws.Range("N29").Value = -6075
$ws.Range("H31").Value = 2192.8572
$ws.Range("I31").Value = 1337.5
$ws.Range("J31").Value = 3333.3333
$ws.Range("K31").Value = 1337.5
$ws.Range("L31").Value = 3333.3333
$ws.Range("M31").Value = -1045.5
$ws.Range("N31").Value = -3917.3333
$ws.Range("H37").Value = 2192.8572
$ws.Range("I37").Value = 1337.5
$ws.Range("J37").Value = 3333.3333
$ws.Range("K37").Value = 1337.5
$ws.Range("L37").Value = 3333.3333
$ws.Range("M37").Value = -1060.5
$ws.Range("N37").Value = -3887.3333
$ws.Range("H99").Value = 10975
$ws.Range("I99").Value = 5910.25
$ws.Range("J99").Value = 21104.5
$ws.Range("K99").Value = 5910.25
$ws.Range("L99").Value = 21104.5
$ws.Range("M99").Value = -3664.25
$ws.Range("N99").Value = -25596.5
$ws.Range("H102").Value = 2314.9312
$ws.Range("I102").Value = 2405.18
$ws.Range("K102").Value = 2405.18
$ws.Range("M102").Value = -783.1799999999998
$ws.Range("H104").Value = 65188.8
$ws.Range("J104").Value = 65188.8
$ws.Range("L104").Value = 65188.8
$ws.Range("N104").Value = -72176.8
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H122").Value = 7266491.5
$ws.Range("I122").Value = 14527473
$ws.Range("J122").Value = 5510.4
$ws.Range("K122").Value = 43582419
$ws.Range("L122").Value = 16531.2
$ws.Range("M122").Value = -43579969
$ws.Range("N122").Value = -21431.2
$ws.Range("H132").Value = 3326.65
$ws.Range("I132").Value = 1214.1818
$ws.Range("J132").Value = 13285.429
$ws.Range("K132").Value = 3642.5454
$ws.Range("L132").Value = 39856.287
$ws.Range("M132").Value = -1112.5454
$ws.Range("N132").Value = -44916.287
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7005.0454
$ws.Range("I7").Value = 5476.2
$ws.Range("J7").Value = 8279.083000000001
$ws.Range("K7").Value = 5476.2
$ws.Range("L7").Value = 8279.083000000001
$ws.Range("M7").Value = -5364.2
$ws.Range("N7").Value = -8503.083000000001
$ws.Range("H23").Value = 5503.5
$ws.Range("J23").Value = 5503.5
$ws.Range("L23").Value = 5503.5
$ws.Range("N23").Value = -5963.5
$ws.Range("H93").Value = 2100.516
$ws.Range("I93").Value = 1243.6842
$ws.Range("J93").Value = 3457.1667
$ws.Range("K93").Value = 1243.6842
$ws.Range("L93").Value = 3457.1667
$ws.Range("M93").Value = 4.315800000000081
$ws.Range("N93").Value = -5953.1667
$ws.Range("H126").Value = 7005.0454
$ws.Range("I126").Value = 5476.2
$ws.Range("J126").Value = 8279.083000000001
$ws.Range("K126").Value = 16428.6
$ws.Range("L126").Value = 24837.249
$ws.Range("M126").Value = -13958.6
$ws.Range("N126").Value = -29777.249
$ws.Range("H132").Value = 6101921
$ws.Range("I132").Value = 10206093
$ws.Range("J132").Value = 7847.727
$ws.Range("K132").Value = 30618279
$ws.Range("L132").Value = 23543.181
$ws.Range("M132").Value = -30615749
$ws.Range("N132").Value = -28603.181
$ws.Range("H136").Value = 6095.0894
$ws.Range("I136").Value = 2094.4146
$ws.Range("J136").Value = 12403.846
$ws.Range("K136").Value = 6283.2438
$ws.Range("L136").Value = 37211.538
$ws.Range("M136").Value = -3733.2438
$ws.Range("N136").Value = -42311.538
$ws.Range("H140").Value = 122337.91
$ws.Range("J140").Value = 122337.91
$ws.Range("L140").Value = 122337.91
$ws.Range("N140").Value = -132697.91
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16396988
$ws.Range("I132").Value = 20411468
$ws.Range("J132").Value = 4529
$ws.Range("K132").Value = 61234404
$ws.Range("L132").Value = 13587
$ws.Range("M132").Value = -61231874
$ws.Range("N132").Value = -18647
$ws.Range("H136").Value = 22958374
$ws.Range("I136").Value = 35715596
$ws.Range("J136").Value = 633232.1
$ws.Range("K136").Value = 107146788
$ws.Range("L136").Value = 1899696.3
$ws.Range("M136").Value = -107144238
$ws.Range("N136").Value = -1904796.3
